$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2057
$ws1.Range("F6").Value = 595
$ws1.Range("F9").Value = 10573
$ws1.Range("F18").Value = 225

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2057
$ws4.Range("F6").Value = 595
$ws4.Range("F12").Value = 10573
$ws4.Range("F21").Value = 225
